$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.48 = 50229.99 pesos`n✅ 50229.99 pesos = 12.42 = 968.83 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 80.13500000000001
$wsTasas.Range("O10").Value = 4025.18

$wsTasas.Range("N12").Value = 4044.5
$wsTasas.Range("O12").Value = 78.01000000000001
